# Apply cell value updates to match the target diff for sheet1 (Jogos_do_Dia_Betfair_Back_Lay_2025-12-22)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F2").Value = 1.01
$ws.Range("G2").Value = 44
$ws.Range("H2").Value = 1.09
$ws.Range("I2").Value = 870
$ws.Range("J2").Value = 1.09
$ws.Range("N2").Value = 1.1
$ws.Range("P2").Value = 1.34
$ws.Range("S2").Value = 1.41
$ws.Range("V2").Value = 1.06
$ws.Range("W2").Value = 1.14
$ws.Range("AD3").Value = 20
$ws.Range("AF3").Value = 11.5
$ws.Range("F3").Value = 1.74
$ws.Range("J3").Value = 4.1
$ws.Range("AB4").Value = 16
$ws.Range("AM4").Value = 270
$ws.Range("AB5").Value = 22
$ws.Range("AE5").Value = 18.5
$ws.Range("AF5").Value = 60
$ws.Range("AG5").Value = 29
$ws.Range("AJ5").Value = 270
$ws.Range("AK5").Value = 130
$ws.Range("AL5").Value = 130
$ws.Range("AN5").Value = 200
$ws.Range("AO5").Value = 9.800000000000001
$ws.Range("G5").Value = 7.8
$ws.Range("H5").Value = 1.54
$ws.Range("J5").Value = 4.2
$ws.Range("L5").Value = 1.4
$ws.Range("N5").Value = 3.7
$ws.Range("O5").Value = 1.32
$ws.Range("P5").Value = 1.94
$ws.Range("Q5").Value = 1.93
$ws.Range("R5").Value = 1.34
$ws.Range("S5").Value = 3.45
$ws.Range("T5").Value = 2
$ws.Range("X5").Value = 19
$ws.Range("J6").Value = 3.25
$ws.Range("AN7").Value = 4.7
$ws.Range("F7").Value = 1.28
$ws.Range("G7").Value = 1.34
$ws.Range("H7").Value = 10
$ws.Range("I7").Value = 13
$ws.Range("J7").Value = 6
$ws.Range("K7").Value = 7.2
$ws.Range("R7").Value = 1.72
$ws.Range("S7").Value = 2.12
$ws.Range("T7").Value = 1.85
$ws.Range("U7").Value = 1.95
$ws.Range("V7").Value = 1.08
$ws.Range("W7").Value = 3.95
$ws.Range("AA8").Value = 90
$ws.Range("AG8").Value = 12.5
$ws.Range("AH8").Value = 21
$ws.Range("AJ8").Value = 29
$ws.Range("AM8").Value = 85
$ws.Range("AN8").Value = 15.5
$ws.Range("G8").Value = 2.14
$ws.Range("N8").Value = 4.3
$ws.Range("O8").Value = 1.24
$ws.Range("P8").Value = 2.08
$ws.Range("Q8").Value = 1.71
$ws.Range("R8").Value = 1.42
$ws.Range("S8").Value = 2.96
$ws.Range("T8").Value = 1.65
$ws.Range("U8").Value = 2.22
$ws.Range("W8").Value = 1.87
$ws.Range("X8").Value = 21
$ws.Range("Y8").Value = 19.5
$ws.Range("AH9").Value = 16.5
$ws.Range("AN9").Value = 34
$ws.Range("G9").Value = 3.65
$ws.Range("J9").Value = 3.45
$ws.Range("Q9").Value = 1.88
$ws.Range("S9").Value = 3.15
$ws.Range("F10").Value = 4
$ws.Range("H10").Value = 1.81
$ws.Range("J10").Value = 3.4
$ws.Range("K10").Value = 4.3
$ws.Range("N10").Value = 4.1
$ws.Range("P10").Value = 2.24
$ws.Range("Q10").Value = 1.64
$ws.Range("R10").Value = 1.5
$ws.Range("S10").Value = 2.56
$ws.Range("T10").Value = 1.62
$ws.Range("U10").Value = 2.28
$ws.Range("AF11").Value = 10
$ws.Range("AH11").Value = 40
$ws.Range("AL11").Value = 42
$ws.Range("AN11").Value = 4.8
$ws.Range("J11").Value = 6.4
$ws.Range("P11").Value = 2.54
$ws.Range("R11").Value = 1.63
$ws.Range("S11").Value = 2.3
$ws.Range("U11").Value = 1.74
$ws.Range("X11").Value = 34
$ws.Range("AA12").Value = 65
$ws.Range("AD12").Value = 16
$ws.Range("AE12").Value = 42
$ws.Range("AO12").Value = 36
$ws.Range("F12").Value = 2.46
$ws.Range("G12").Value = 2.68
$ws.Range("H12").Value = 2.86
$ws.Range("I12").Value = 3.15
$ws.Range("J12").Value = 3.4
$ws.Range("K12").Value = 3.7
$ws.Range("N12").Value = 3.75
$ws.Range("O12").Value = 1.31
$ws.Range("Q12").Value = 1.91
$ws.Range("R12").Value = 1.37
$ws.Range("U12").Value = 2.18
$ws.Range("V12").Value = 1.46
$ws.Range("W12").Value = 1.6
$ws.Range("Y12").Value = 15
$ws.Range("Z12").Value = 26
$ws.Range("G13").Value = 1.71
$ws.Range("P13").Value = 1.59
$ws.Range("S13").Value = 3.05
$ws.Range("W13").Value = 2.42
$ws.Range("AE14").Value = 34
$ws.Range("AM14").Value = 70
$ws.Range("F14").Value = 2.34
$ws.Range("G14").Value = 2.56
$ws.Range("H14").Value = 2.84
$ws.Range("J14").Value = 3.8
$ws.Range("Q14").Value = 1.57
$ws.Range("S14").Value = 2.4
$ws.Range("T14").Value = 1.52
$ws.Range("U14").Value = 2.52
$ws.Range("W14").Value = 1.68
$ws.Range("X14").Value = 26
$ws.Range("Z14").Value = 28
$ws.Range("AA15").Value = 10
$ws.Range("AB15").Value = 32
$ws.Range("AE15").Value = 16
$ws.Range("AF15").Value = 120
$ws.Range("AJ15").Value = 630
$ws.Range("AK15").Value = 280
$ws.Range("AL15").Value = 230
$ws.Range("AN15").Value = 440
$ws.Range("AO15").Value = 7
$ws.Range("F15").Value = 12
$ws.Range("G15").Value = 12.5
$ws.Range("H15").Value = 1.37
$ws.Range("I15").Value = 1.38
$ws.Range("J15").Value = 5.3
$ws.Range("K15").Value = 5.4
$ws.Range("N15").Value = 3.9
$ws.Range("O15").Value = 1.32
$ws.Range("P15").Value = 1.98
$ws.Range("Q15").Value = 1.96
$ws.Range("T15").Value = 2.44
$ws.Range("U15").Value = 1.62
$ws.Range("V15").Value = 3.65
$ws.Range("W15").Value = 1.08
$ws.Range("X15").Value = 15.5
$ws.Range("Y15").Value = 7.2
$ws.Range("Z15").Value = 7
$ws.Range("AG16").Value = 20
$ws.Range("AN16").Value = 34
$ws.Range("I16").Value = 1.82
$ws.Range("J16").Value = 4.3
$ws.Range("K16").Value = 4.9
$ws.Range("S16").Value = 2.18
$ws.Range("V16").Value = 2.22
$ws.Range("AB17").Value = 16.5
$ws.Range("AC17").Value = 10.5
$ws.Range("AG17").Value = 13.5
$ws.Range("AH17").Value = 18
$ws.Range("AL17").Value = 34
$ws.Range("AM17").Value = 70
$ws.Range("AN17").Value = 17
$ws.Range("AO17").Value = 18
$ws.Range("F17").Value = 2.58
$ws.Range("G17").Value = 2.74
$ws.Range("H17").Value = 2.52
$ws.Range("I17").Value = 2.68
$ws.Range("J17").Value = 3.95
$ws.Range("K17").Value = 4.3
$ws.Range("S17").Value = 2.5
$ws.Range("T17").Value = 1.56
$ws.Range("U17").Value = 2.54
$ws.Range("V17").Value = 1.59
$ws.Range("W17").Value = 1.57
$ws.Range("X17").Value = 26
$ws.Range("Y17").Value = 18
$ws.Range("Z17").Value = 22
$ws.Range("M18").Value = 1.09
$ws.Range("Q18").Value = 2.2
$ws.Range("S18").Value = 4.2
$ws.Range("U18").Value = 1.74
$ws.Range("V18").Value = 1.16
$ws.Range("AA19").Value = 55
$ws.Range("AM19").Value = 95
$ws.Range("F19").Value = 2.58
$ws.Range("G19").Value = 2.6
$ws.Range("H19").Value = 3.1
$ws.Range("P19").Value = 1.85
$ws.Range("Q19").Value = 2.14
$ws.Range("V19").Value = 1.46
$ws.Range("Z19").Value = 19.5
$ws.Range("AB20").Value = 7
$ws.Range("AE20").Value = 95
$ws.Range("AM20").Value = 190
$ws.Range("F20").Value = 1.86
$ws.Range("G20").Value = 1.87
$ws.Range("J20").Value = 3.55
$ws.Range("K20").Value = 3.6
$ws.Range("R20").Value = 1.24
$ws.Range("W20").Value = 2.14
$ws.Range("AN21").Value = 5.8
$ws.Range("F21").Value = 1.33
$ws.Range("H21").Value = 11
$ws.Range("K21").Value = 6.2
$ws.Range("P21").Value = 2.22
$ws.Range("Q21").Value = 1.75
$ws.Range("S21").Value = 2.92
